$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 263.1
$ws.Cells.Item(8, 9).Value = 263.1
$ws.Cells.Item(8, 11).Value = 789.3000000000001
$ws.Cells.Item(8, 13).Value = -650.3000000000001
$ws.Cells.Item(62, 14).Value = -6860.5
$ws.Cells.Item(62, 11).Value = 2456.111
$ws.Cells.Item(62, 13).Value = -1832.111
$ws.Cells.Item(62, 8).Value = 3718.6667
$ws.Cells.Item(62, 12).Value = 5612.5
$ws.Cells.Item(62, 9).Value = 2456.111
$ws.Cells.Item(62, 10).Value = 5612.5
$ws.Cells.Item(65, 13).Value = -9160.555
$ws.Cells.Item(65, 9).Value = 2456.111
$ws.Cells.Item(65, 14).Value = -34302.5
$ws.Cells.Item(65, 10).Value = 5612.5
$ws.Cells.Item(65, 11).Value = 12280.555
$ws.Cells.Item(65, 8).Value = 3718.6667
$ws.Cells.Item(65, 12).Value = 28062.5
$ws.Cells.Item(98, 8).Value = 4622.222
$ws.Cells.Item(98, 13).Value = -3124.222
$ws.Cells.Item(98, 9).Value = 4622.222
$ws.Cells.Item(98, 11).Value = 4622.222
$ws.Cells.Item(122, 9).Value = 4622.222
$ws.Cells.Item(122, 13).Value = -11416.666
$ws.Cells.Item(122, 8).Value = 4622.222
$ws.Cells.Item(122, 11).Value = 13866.666
$ws.Cells.Item(125, 10).Value = 879.6
$ws.Cells.Item(125, 12).Value = 7916.400000000001
$ws.Cells.Item(125, 14).Value = -12836.4
$ws.Cells.Item(125, 8).Value = 879.6
$ws.Cells.Item(137, 14).Value = -7158092.399999999
$ws.Cells.Item(137, 13).Value = -7693.7145
$ws.Cells.Item(137, 10).Value = 2384330.8
$ws.Cells.Item(137, 8).Value = 1193872.6
$ws.Cells.Item(137, 11).Value = 10243.7145
$ws.Cells.Item(137, 9).Value = 3414.5715
$ws.Cells.Item(137, 12).Value = 7152992.399999999
$ws.Cells.Item(141, 8).Value = 4581.364
$ws.Cells.Item(141, 13).Value = -5843.125
$ws.Cells.Item(141, 9).Value = 3674.375
$ws.Cells.Item(141, 11).Value = 11023.125
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 14).Value = -4230.8
$ws.Cells.Item(32, 11).Value = 23720.936
$ws.Cells.Item(32, 9).Value = 23720.936
$ws.Cells.Item(32, 13).Value = -23433.936
$ws.Cells.Item(32, 10).Value = 3656.8
$ws.Cells.Item(32, 12).Value = 3656.8
$ws.Cells.Item(32, 8).Value = 18787.13
$ws.Cells.Item(61, 12).Value = 5664.6665
$ws.Cells.Item(61, 10).Value = 5664.6665
$ws.Cells.Item(61, 8).Value = 4821.5947
$ws.Cells.Item(61, 14).Value = -6088.6665
$ws.Cells.Item(61, 9).Value = 4658.4194
$ws.Cells.Item(61, 11).Value = 4658.4194
$ws.Cells.Item(61, 13).Value = -4446.4194
$ws.Cells.Item(132, 11).Value = 5640.916499999999
$ws.Cells.Item(132, 9).Value = 1880.3055
$ws.Cells.Item(132, 14).Value = -13995.5
$ws.Cells.Item(132, 12).Value = 8935.5
$ws.Cells.Item(132, 8).Value = 1990.125
$ws.Cells.Item(132, 13).Value = -3110.916499999999
$ws.Cells.Item(132, 10).Value = 2978.5
$ws.Cells.Item(136, 8).Value = 4821.5947
$ws.Cells.Item(136, 13).Value = -11425.2582
$ws.Cells.Item(136, 9).Value = 4658.4194
$ws.Cells.Item(136, 11).Value = 13975.2582
$ws.Cells.Item(136, 10).Value = 5664.6665
$ws.Cells.Item(136, 12).Value = 16993.9995
$ws.Cells.Item(136, 14).Value = -22093.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 14).Value = -56226
$ws.Cells.Item(2, 8).Value = 56000
$ws.Cells.Item(2, 12).Value = 56000
$ws.Cells.Item(2, 10).Value = 56000
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 14).ClearContents() | Out-Null
$ws.Cells.Item(134, 8).Value = 4391.5
$ws.Cells.Item(134, 12).Value = 23910
$ws.Cells.Item(134, 10).Value = 7970
$ws.Cells.Item(134, 14).Value = -28980
$ws.Cells.Item(134, 13).Value = -8850.249899999999
$ws.Cells.Item(134, 11).Value = 11385.2499
$ws.Cells.Item(134, 9).Value = 3795.0833
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 11).Value = 5678.095
$ws.Cells.Item(31, 8).Value = 826981.5600000001
$ws.Cells.Item(31, 13).Value = -5383.095
$ws.Cells.Item(31, 9).Value = 5678.095
$ws.Cells.Item(34, 13).Value = -5476.095
$ws.Cells.Item(34, 11).Value = 5678.095
$ws.Cells.Item(34, 8).Value = 826981.5600000001
$ws.Cells.Item(34, 9).Value = 5678.095
$ws.Cells.Item(99, 11).Value = 1072.5
$ws.Cells.Item(99, 8).Value = 3858
$ws.Cells.Item(99, 13).Value = 425.5
$ws.Cells.Item(99, 9).Value = 1072.5
$ws.Cells.Item(99, 14).Value = -17996
$ws.Cells.Item(99, 10).Value = 15000
$ws.Cells.Item(99, 12).Value = 15000
$ws.Cells.Item(126, 8).Value = 3858
$ws.Cells.Item(126, 14).Value = -49940
$ws.Cells.Item(126, 11).Value = 3217.5
$ws.Cells.Item(126, 13).Value = -747.5
$ws.Cells.Item(126, 12).Value = 45000
$ws.Cells.Item(126, 10).Value = 15000
$ws.Cells.Item(126, 9).Value = 1072.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 12).Value = 13334.625
$ws.Cells.Item(68, 10).Value = 4444.875
$ws.Cells.Item(68, 8).Value = 2874.0688
$ws.Cells.Item(68, 14).Value = -14956.625
$ws.Cells.Item(71, 10).Value = 4444.875
$ws.Cells.Item(71, 14).Value = -48115.875
$ws.Cells.Item(71, 12).Value = 40003.875
$ws.Cells.Item(71, 8).Value = 2874.0688
$ws.Cells.Item(107, 11).Value = 864.5581199999999
$ws.Cells.Item(107, 10).Value = 2072.9473
$ws.Cells.Item(107, 9).Value = 288.18604
$ws.Cells.Item(107, 13).Value = 1055.44188
$ws.Cells.Item(107, 12).Value = 6218.841899999999
$ws.Cells.Item(107, 14).Value = -10058.8419
$ws.Cells.Item(107, 8).Value = 835.129
$ws.Cells.Item(122, 8).Value = 1249.5
$ws.Cells.Item(122, 10).Value = 1499
$ws.Cells.Item(122, 12).Value = 13491
$ws.Cells.Item(122, 14).Value = -18391
$ws.Cells.Item(130, 14).Value = -33536
$ws.Cells.Item(130, 10).Value = 7832
$ws.Cells.Item(130, 11).Value = 9090
$ws.Cells.Item(130, 12).Value = 23496
$ws.Cells.Item(130, 8).Value = 5431
$ws.Cells.Item(130, 9).Value = 3030
$ws.Cells.Item(130, 13).Value = -4070
$ws.Cells.Item(133, 13).Value = -2603.3329
$ws.Cells.Item(133, 10).Value = 5749.5835
$ws.Cells.Item(133, 9).Value = 2554.4443
$ws.Cells.Item(133, 8).Value = 4380.2383
$ws.Cells.Item(133, 12).Value = 17248.7505
$ws.Cells.Item(133, 14).Value = -27368.7505
$ws.Cells.Item(133, 11).Value = 7663.3329
$ws.Cells.Item(134, 8).Value = 5132.273
$ws.Cells.Item(134, 12).Value = 12653.571
$ws.Cells.Item(134, 10).Value = 4217.857
$ws.Cells.Item(134, 14).Value = -22793.571
$ws.Cells.Item(134, 13).Value = -15127.5
$ws.Cells.Item(134, 11).Value = 20197.5
$ws.Cells.Item(134, 9).Value = 6732.5
$ws.Cells.Item(137, 14).Value = -395274.75
$ws.Cells.Item(137, 13).Value = 1506
$ws.Cells.Item(137, 10).Value = 128358.25
$ws.Cells.Item(137, 8).Value = 102926.2
$ws.Cells.Item(137, 11).Value = 3594
$ws.Cells.Item(137, 9).Value = 1198
$ws.Cells.Item(137, 12).Value = 385074.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 11).Value = 5766.6665
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 8).Value = 5766.6665
$ws.Cells.Item(31, 13).Value = -5474.6665
$ws.Cells.Item(31, 9).Value = 5766.6665
$ws.Cells.Item(31, 14).ClearContents() | Out-Null
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 9).Value = 5766.6665
$ws.Cells.Item(37, 14).ClearContents() | Out-Null
$ws.Cells.Item(37, 8).Value = 5766.6665
$ws.Cells.Item(37, 13).Value = -5489.6665
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 5766.6665
$ws.Cells.Item(55, 11).Value = 4225
$ws.Cells.Item(55, 13).Value = -3898
$ws.Cells.Item(55, 9).Value = 4225
$ws.Cells.Item(55, 8).Value = 4225
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).ClearContents() | Out-Null
$ws.Cells.Item(82, 8).Value = 31721.166
$ws.Cells.Item(82, 14).Value = -32487.166
$ws.Cells.Item(82, 12).Value = 31721.166
$ws.Cells.Item(82, 10).Value = 31721.166
$ws.Cells.Item(85, 8).Value = 31721.166
$ws.Cells.Item(85, 12).Value = 31721.166
$ws.Cells.Item(85, 10).Value = 31721.166
$ws.Cells.Item(85, 14).Value = -34373.166
$ws.Cells.Item(122, 9).Value = 16911.445
$ws.Cells.Item(122, 13).Value = -48284.335
$ws.Cells.Item(122, 8).Value = 15720.3
$ws.Cells.Item(122, 11).Value = 50734.335
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(126, 8).Value = 2448.5715
$ws.Cells.Item(126, 14).Value = -16040
$ws.Cells.Item(126, 11).Value = 5260.0002
$ws.Cells.Item(126, 13).Value = -2790.0002
$ws.Cells.Item(126, 12).Value = 11100
$ws.Cells.Item(126, 10).Value = 3700
$ws.Cells.Item(126, 9).Value = 1753.3334
$ws.Cells.Item(132, 14).Value = -11058.5
$ws.Cells.Item(132, 12).Value = 5998.5
$ws.Cells.Item(132, 8).Value = 1557.1333
$ws.Cells.Item(132, 10).Value = 1999.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 8).Value = 2800.8
$ws.Cells.Item(7, 14).Value = -5224
$ws.Cells.Item(40, 10).Value = 4633.3335
$ws.Cells.Item(40, 8).Value = 3800
$ws.Cells.Item(40, 12).Value = 4633.3335
$ws.Cells.Item(40, 14).Value = -4905.3335
$ws.Cells.Item(40, 9).Value = 3487.5
$ws.Cells.Item(40, 11).Value = 3487.5
$ws.Cells.Item(40, 13).Value = -3351.5
$ws.Cells.Item(117, 12).Value = 62900
$ws.Cells.Item(117, 14).Value = -72078
$ws.Cells.Item(117, 10).Value = 62900
$ws.Cells.Item(117, 8).Value = 62900
$ws.Cells.Item(122, 9).Value = 6531.1304
$ws.Cells.Item(122, 13).Value = -17143.3912
$ws.Cells.Item(122, 8).Value = 6210.8965
$ws.Cells.Item(122, 11).Value = 19593.3912
$ws.Cells.Item(122, 10).Value = 4983.3335
$ws.Cells.Item(122, 12).Value = 14950.0005
$ws.Cells.Item(122, 14).Value = -19850.0005
$ws.Cells.Item(126, 8).Value = 2800.8
$ws.Cells.Item(126, 14).Value = -19940
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 5476.9998
$ws.Cells.Item(132, 9).Value = 1825.6666
$ws.Cells.Item(132, 14).Value = -14902.0001
$ws.Cells.Item(132, 12).Value = 9842.000100000001
$ws.Cells.Item(132, 8).Value = 2553.1667
$ws.Cells.Item(132, 13).Value = -2946.9998
$ws.Cells.Item(132, 10).Value = 3280.6667
$ws.Cells.Item(134, 8).Value = 58499.332
$ws.Cells.Item(134, 12).Value = 58499.332
$ws.Cells.Item(134, 10).Value = 58499.332
$ws.Cells.Item(134, 14).Value = -68639.33199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1912.75
$ws.Cells.Item(126, 14).Value = -11237
$ws.Cells.Item(126, 11).Value = 5339.142599999999
$ws.Cells.Item(126, 13).Value = -2869.142599999999
$ws.Cells.Item(126, 12).Value = 6297
$ws.Cells.Item(126, 10).Value = 2099
$ws.Cells.Item(126, 9).Value = 1779.7142
$ws.Cells.Item(136, 8).Value = 5167.511
$ws.Cells.Item(136, 13).Value = -2392.7142
$ws.Cells.Item(136, 9).Value = 1647.5714
$ws.Cells.Item(136, 11).Value = 4942.7142
$ws.Cells.Item(136, 10).Value = 8247.458000000001
$ws.Cells.Item(136, 12).Value = 24742.374
$ws.Cells.Item(136, 14).Value = -29842.374

Write-Host "Applied all changes"